# Convert the old "config" / "sapLogin" prototype into the new
# "Rutas" / "parametrosInicio" layout.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # was "config"   -> "Rutas"
$ws2 = $wb.Worksheets.Item(2)   # was "sapLogin" -> "parametrosInicio"

# --- Rename the sheets -------------------------------------------------
$ws1.Name = "Rutas"
$ws2.Name = "parametrosInicio"

# --- Sheet 1 ("Rutas"): drop every row except the SAP-logon path row ----
# (keep row 2, just relabel A2, clear the rest of the old rows)
$ws1.Range("A1:B1").ClearContents()
$ws1.Range("A2").Value = "Path SAP logon"
$ws1.Range("A3:B3").ClearContents()
$ws1.Range("A4:B4").ClearContents()
$ws1.Range("A5:B5").ClearContents()

# leftover cursor position from the author's last session
[void]$ws1.Range("B10").Select()

# --- Sheet 2 ("parametrosInicio"): keep usuario/psw, add the new rows --
# Row 3: new "ambiente" parameter (its value is later driven by the
# validation dropdown below)
$ws2.Range("A3").Value = "ambiente"
$ws2.Range("B3").Value = "QAS - EHP8 on HANA"

# Row 5 / Row 6: new "Fecha" / "Periodo" parameters.
# Intern "Periodo" before "30.10.2022" so the shared-string table ends up
# in the same order the original author produced.
$ws2.Range("A5").Value = "Fecha"
$ws2.Range("A6").Value = "Periodo"
$ws2.Range("B5").Value = "30.10.2022"
$ws2.Range("B6").Value = 7

# Column widths now that "QAS - EHP8 on HANA" / "ambiente" live here
$ws2.Columns.Item(1).ColumnWidth = 11.65
$ws2.Columns.Item(2).ColumnWidth = 18.5

# Dropdown list validation on the "ambiente" value cell
[void]$ws2.Range("B3").Validation.Add(3, 1, 1, '"QAS - EHP8 on HANA,PRD SAP HANA"')

[void]$ws2.Range("B6").Select()

# "parametrosInicio" is now the front-most / active tab
[void]$ws2.Activate()
